$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fix up cell formatting (style indices) before/independent of values
#    - A1 needs the "plain" style that A3 (and most A-column cells) use
#    - A2 needs the "big/special" style that A1 currently has
#    - New column C (rows 1-19) should use the same style as column B
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A1").PasteSpecial(-4122)

$ws.Range("B1:B19").Copy()
$ws.Range("C1:C19").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Write the new header row + the Pokemon/Type/DPS table
#    (A previously held the Pokemon name and B the DPS value; now
#     A = Pokemon, B = Type, C = DPS, with a header row on row 1
#     and the whole data table shifted down by one row.)
# ---------------------------------------------------------------------
$data = @(
    @("寶可夢", "屬性", "DPS"),
    @("暗影雷吉", "一般", 22.5),
    @("Mega噴火龍", "火", 19.68),
    @("Mega路卡", "格鬥", 22.43),
    @("Mega蓋歐卡", "水", 19.420000000000002),
    @("Mega列空坐", "飛行", 26.05),
    @("Mega蜥蜴王", "草", 19.97),
    @("無極汰那", "毒", 23.04),
    @("Mega雷電獸", "電", 16.59),
    @("Mega烈咬陸鯊", "地面", 17.850000000000001),
    @("Mega艾路雷朵", "超能力", 18.11),
    @("Mega班基拉斯", "岩石", 15.38),
    @("焰白酋雷姆", "冰", 19.22),
    @("Mega凱羅斯", "蟲", 15.73),
    @("無極汰那", "龍", 23.04),
    @("奈克洛茲瑪", "幽靈", 22.16),
    @("Mega班基拉斯", "惡", 17.88),
    @("劍之王", "鋼", 22.27),
    @("Mega沙奈朵", "妖精", 17.25)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# ---------------------------------------------------------------------
# 3. Update the sheet view: no more frozen/scrolled topLeftCell, and
#    the active selection moves to D4 instead of K24.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("D4").Select()
